$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns B and C
$ws.Range("B1").Value = "值"
$ws.Range("C1").Value = "出處"

# New rows further down the sheet
$ws.Range("A13").Value = "備註"
$ws.Range("A15").Value = "衝突"

# Set column widths for the newly used columns B and C
# (values chosen so the stored OOXML width matches the target as closely as
# this runtime's character-width rounding allows: target widths are 12 and 52.796875)
$ws.Columns.Item(2).ColumnWidth = 11.285714285714286
$ws.Columns.Item(3).ColumnWidth = 52.142857142857146

# Update the active selection to match the target state
$ws.Range("C13").Select()
